$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 (shifts existing rows 59:67 down to 60:68)
$ws.Rows.Item(59).Insert()

# The new row 59 is a copy of the former row 59 (now row 60) with a handful
# of values changed: Fecha, Precio mínimo, Precio máximo, Precio promedio
# ponderado, Precio $/Kg.
$ws.Cells.Item(59, 1).Value = 2
$ws.Cells.Item(59, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44889
$ws.Cells.Item(59, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = 100112032
$ws.Cells.Item(59, 7).Value = "Zapallo italiano"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 500
$ws.Cells.Item(59, 11).Value = 6000
$ws.Cells.Item(59, 12).Value = 7000
$ws.Cells.Item(59, 13).Value = 6500
$ws.Cells.Item(59, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(59, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(59, 16).Value = 108
$ws.Cells.Item(59, 17).Value = 60
$ws.Cells.Item(59, 18).Value = "Hortaliza"
